$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B49 from text "4" to numeric 4
$ws.Range("B49").Value = 4

# Add new row 50
$ws.Range("A50").Value = "Ruilin"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "3"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "无"
$ws.Range("D50").Value = "FBK"
$ws.Range("E50").Value = "OTH"
$ws.Range("F50").Value = "6b52664d-9c16-43d6-8192-f8f8d8ec6227"
$ws.Range("G50").Value = "SkBYYyZRZ_annotated.xlsx"
$ws.Range("H50").Value = "We agree, and we will open source some of the experiments around the time of acceptance."
